$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-19 from 45184 to 45185
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
